$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: G3 and H3 -> 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: D4 and E4 -> 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: D5 and E5 -> 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Rows 6-18: H column -> 1
for ($r = 6; $r -le 18; $r++) {
    $ws.Range("H$r").Value = 1
}
